# Update generator data to NETL data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: ratio values used by the new formulas in B2/E2
$ws.Range("A5").Value = "ratio from ES-4"
$ws.Range("B5").Value = 1.49
$ws.Range("E5").Value = 2.39

# Row 2 data updates (NETL data)
$ws.Range("B2").Formula = "=C2*B5"
$ws.Range("C2").Value = 64000
$ws.Range("E2").Formula = "=F2*E5"
$ws.Range("F2").Value = 13000
$ws.Range("G2").Value = 10000
$ws.Range("H2").Value = 124000
$ws.Range("I2").Value = 16000
$ws.Range("J2").Value = 38000

# Column A width to fit the new "ratio from ES-4" label
$ws.Range("A1").EntireColumn.ColumnWidth = 12.83

# Selection matches the authored state
$ws.Range("E3").Select() | Out-Null
